$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2:O30").Value = "2022-07-23 20:59:31"
